$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "hpi"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("hpi")

$ws.Range("D2").Value = "Absence of relief from antacids suggests that the pain may not be related to acid reflux, which is less common in achalasia."

$ws.Range("B4").Value = "Dysphagia is a common symptom associated with achalasia, which is relevant to the diagnosis."
$ws.Range("D4").Value = "Coughing can be associated with esophageal issues; its absence may suggest a lower likelihood of esophageal complications."

$ws.Range("B5").Value = "Dysphagia for solids further supports the presence of achalasia, which is consistent with the diagnosis."
$ws.Range("D5").Value = "Nausea and vomiting are common in esophageal disorders; their absence may indicate a different underlying issue."

$ws.Range("D6").Value = "Shortness of breath can indicate pulmonary complications or severe esophageal issues, which are not present here."

# ---------------------------------------------------------------------------
# Sheet "hist"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("hist")

$ws.Range("B2").Value = "CREST syndrome is more prevalent in females, which supports the likelihood of having the condition."
$ws.Range("D2").Value = "Radiation exposure is a known risk factor for esophageal conditions; its absence may argue against achalasia."

$ws.Range("B3").Value = "CREST syndrome typically presents in middle-aged individuals, making this finding supportive of the diagnosis."
$ws.Range("C3").Value = "Absence of diagnosed hypertension"
$ws.Range("D3").Value = "Hypertension is often associated with vascular conditions; its absence may reduce the likelihood of CREST syndrome."

$ws.Range("C4").Value = "Absence of prior myocardial infarction"
$ws.Range("D4").Value = "A history of myocardial infarction can indicate cardiovascular issues related to CREST syndrome; its absence may suggest lower risk."

$ws.Range("C5").Value = "Absence of type 2 diabetes"
$ws.Range("D5").Value = "Type 2 diabetes is a risk factor for vascular complications; its absence may argue against the presence of CREST syndrome."

$ws.Range("B6").Value = "Nicotine dependence can exacerbate vascular issues; its absence may favor a diagnosis of CREST syndrome."
$ws.Range("C6").Value = "Absence of obesity"
$ws.Range("D6").Value = "Obesity is a risk factor for various conditions, including achalasia; its absence may suggest a lower likelihood of the diagnosis."

# ---------------------------------------------------------------------------
# Sheet "soc"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("soc")

$ws.Range("D3").Value = "The absence of a family history of cancer may suggest a lower likelihood of genetic predispositions that could be linked to CREST syndrome."
$ws.Range("D4").Value = "Absence of recent travel may indicate a lack of exposure to infections or environmental factors that could trigger or exacerbate autoimmune conditions."
$ws.Range("D5").Value = "The absence of recent medical procedures may suggest a lack of acute exacerbations or complications that could be associated with CREST syndrome."
$ws.Range("D6").Value = "The absence of gestational complications may indicate a lower likelihood of underlying autoimmune conditions that could affect pregnancy, such as CREST syndrome."

# ---------------------------------------------------------------------------
# Sheet "obj"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("obj")

$ws.Range("C2").Value = "Hand thickening observed is absent."
$ws.Range("D2").Value = "Hand thickening is a characteristic feature of CREST syndrome; its absence strongly suggests that the syndrome is not present."

$ws.Range("B3").Value = "Weight loss can occur in patients with achalasia due to difficulty swallowing and subsequent malnutrition, which may support the diagnosis of Type 2 Achalasia."
$ws.Range("C3").Value = "Finger ulcers observed is absent."
$ws.Range("D3").Value = "Finger ulcers are another common manifestation of CREST syndrome; their absence further supports the likelihood of not having the syndrome."

$ws.Range("A4").Value = "Raynaud's phenomenon on exam is absent."
$ws.Range("B4").Value = "While the absence of Raynaud's phenomenon does not confirm the diagnosis, it does not rule it out, as not all patients with CREST syndrome exhibit this finding."
$ws.Range("C4").Value = "Weakness on exam is absent."
$ws.Range("D4").Value = "Weakness is often associated with systemic involvement in CREST syndrome; its absence may indicate that the syndrome is not present."

$ws.Range("A5").Value = "Hoarse voice observed is absent."
$ws.Range("B5").Value = "The absence of a hoarse voice does not negate the diagnosis of achalasia, as it is not a definitive symptom of the condition."
$ws.Range("C5").Value = "Obesity by vital signs is absent."
$ws.Range("D5").Value = "While obesity is not directly related to CREST syndrome, its absence may suggest a lack of systemic involvement that could be expected in such patients."

$ws.Range("A6").Value = "Cough observed is absent."
$ws.Range("B6").Value = "The absence of cough does not provide strong evidence against achalasia, as it is not a primary symptom of the condition."
$ws.Range("C6").Value = "High blood pressure when checked is absent."
$ws.Range("D6").Value = "The absence of high blood pressure does not support the diagnosis of CREST syndrome, as hypertension can be a common finding in other conditions."
